# Horarios actualizados Linea 141 - 951
# Refresh the scraped-schedule workbook: new "Ultima actualizacion" / "Total
# filas" header values on each sheet, and a refreshed data table (some rows
# keep their previous Hora_Scrap stamp, several get new arrival estimates,
# and a batch of new rows is appended) for sheets LP1912, LP1912-215 and
# 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:19:59"
$ws1.Range("A3").Value = "Total filas: 36"

$flat1 = @(
    '05:44:02','05:44','14_ABASTO',0,'LP1912',
    '05:44:02','05:47','17_ROMERO',3,'LP1912',
    '05:58:04','06:09','10_OLMOS',11,'LP1912',
    '05:58:04','06:16','215A_EL PATO',18,'LP1912',
    '06:19:59','06:29','23_HERNANDEZ',10,'LP1912',
    '06:19:59','06:30','23_HERNANDEZ',32,'LP1912',
    '06:19:59','06:33','11_ETCHEVERRY',14,'LP1912',
    '05:58:04','06:34','11_ETCHEVERRY',36,'LP1912',
    '06:19:59','06:38','17X38_ROMERO',19,'LP1912',
    '06:19:59','06:39','17X38_ROMERO',41,'LP1912',
    '05:44:02','06:40','17X38_ROMERO',56,'LP1912',
    '05:58:04','06:41','16_SANTA ANA',43,'LP1912',
    '06:19:59','06:56','215A_EL PATO',37,'LP1912',
    '05:58:04','06:57','215A_EL PATO',59,'LP1912',
    '06:19:59','06:58','225_GOMEZ',39,'LP1912',
    '05:58:04','06:59','225_GOMEZ',61,'LP1912',
    '06:19:59','07:15','215C_EL PATO',56,'LP1912',
    '05:58:04','07:16','215C_EL PATO',78,'LP1912',
    '06:19:59','07:18','14_ABASTO',59,'LP1912',
    '05:58:04','07:19','14_ABASTO',81,'LP1912',
    '06:19:59','07:20','16_SANTA ANA',61,'LP1912',
    '06:19:59','07:21','23_HERNANDEZ',62,'LP1912',
    '05:58:04','07:21','16_SANTA ANA',83,'LP1912',
    '05:58:04','07:22','23_HERNANDEZ',84,'LP1912',
    '06:19:59','07:29','17X38_ROMERO',70,'LP1912',
    '06:19:59','07:34','10_OLMOS',75,'LP1912',
    '05:58:04','07:35','10_OLMOS',97,'LP1912',
    '06:19:59','07:36','27_EL RETIRO',77,'LP1912',
    '05:58:04','07:37','27_EL RETIRO',99,'LP1912',
    '06:19:59','07:43','215A_EL PATO',84,'LP1912',
    '06:19:59','07:54','14_ABASTO',95,'LP1912',
    '05:58:04','07:55','14_ABASTO',117,'LP1912',
    '06:19:59','07:59','17_ROMERO',100,'LP1912',
    '06:19:59','08:00','16_SANTA ANA',101,'LP1912',
    '06:19:59','08:11','10_OLMOS',112,'LP1912',
    '06:19:59','08:12','15X38_ABASTO',113,'LP1912'
)

$n1 = $flat1.Length / 5
for ($i = 0; $i -lt $n1; $i++) {
    $r = 6 + $i
    $base = $i * 5
    $ws1.Cells.Item($r, 1).Value = $flat1[$base]
    $ws1.Cells.Item($r, 2).Value = $flat1[$base + 1]
    $ws1.Cells.Item($r, 3).Value = $flat1[$base + 2]
    $ws1.Cells.Item($r, 4).Value = $flat1[$base + 3]
    $ws1.Cells.Item($r, 5).Value = $flat1[$base + 4]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:19:59"
$ws2.Range("A3").Value = "Total filas: 6"

$flat2 = @(
    '05:58:04','06:16','215A_EL PATO',18,'LP1912',
    '06:19:59','06:56','215A_EL PATO',37,'LP1912',
    '05:58:04','06:57','215A_EL PATO',59,'LP1912',
    '06:19:59','07:15','215C_EL PATO',56,'LP1912',
    '05:58:04','07:16','215C_EL PATO',78,'LP1912',
    '06:19:59','07:43','215A_EL PATO',84,'LP1912'
)

$n2 = $flat2.Length / 5
for ($i = 0; $i -lt $n2; $i++) {
    $r = 6 + $i
    $base = $i * 5
    $ws2.Cells.Item($r, 1).Value = $flat2[$base]
    $ws2.Cells.Item($r, 2).Value = $flat2[$base + 1]
    $ws2.Cells.Item($r, 3).Value = $flat2[$base + 2]
    $ws2.Cells.Item($r, 4).Value = $flat2[$base + 3]
    $ws2.Cells.Item($r, 5).Value = $flat2[$base + 4]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 06:19:59"
$ws3.Range("A3").Value = "Total filas: 2"

$flat3 = @(
    '06:19:59','07:42','215A_LA PLATA',83,'L6173',
    '05:58:04','07:43','215A_LA PLATA',105,'L6173'
)

$n3 = $flat3.Length / 5
for ($i = 0; $i -lt $n3; $i++) {
    $r = 6 + $i
    $base = $i * 5
    $ws3.Cells.Item($r, 1).Value = $flat3[$base]
    $ws3.Cells.Item($r, 2).Value = $flat3[$base + 1]
    $ws3.Cells.Item($r, 3).Value = $flat3[$base + 2]
    $ws3.Cells.Item($r, 4).Value = $flat3[$base + 3]
    $ws3.Cells.Item($r, 5).Value = $flat3[$base + 4]
}
